$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = '''44.569.45'
$c.Style = "Normal"
$ws.Range('E2').Value = '  +3.47%  '
$c = $ws.Range('D3')
$c.Value = '''2.433.80'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +2.29%  '
$c = $ws.Range('D4')
$c.Value = '''1.00'
$c.Style = "Normal"
$ws.Range('E4').Value = '  -0.02%  '
$c = $ws.Range('D5')
$c.Value = '''312.29'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +3.40%  '
$c = $ws.Range('D6')
$c.Value = '''101.87'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +5.11%  '
$ws.Range('E7').Value = '  +1.65%  '
$ws.Range('E8').Value = '  -0.03%  '
$c = $ws.Range('D9')
$c.Value = '''0.508'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +1.52%  '
$c = $ws.Range('D10')
$c.Value = '''35.37'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +3.06%  '
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('E12').Value = '  +1.02%  '
$c = $ws.Range('D13')
$c.Value = '''18.78'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +2.92%  '
$c = $ws.Range('D14')
$c.Value = '''6.97'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +2.39%  '
$c = $ws.Range('D15')
$c.Value = '''2.811.53'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +2.32%  '
$c = $ws.Range('D16')
$c.Value = '''2.453.04'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +3.02%  '
$c = $ws.Range('D17')
$c.Value = '''0.838'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +4.04%  '
$c = $ws.Range('D18')
$c.Value = '''44.493.87'
$c.Style = "Normal"
$c = $ws.Range('D19')
$c.Value = '''12.49'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +2.31%  '
$ws.Range('E20').Value = '  +1.54%  '
$c = $ws.Range('D21')
$c.Value = '''0.0₃0909'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +2.34%  '
$c = $ws.Range('D22')
$c.Value = '''68.95'
$c.Style = "Normal"
$c = $ws.Range('D23')
$c.Value = '''2.32'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +3.23%  '
$c = $ws.Range('D24')
$c.Value = '''241.14'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +2.38%  '
$ws.Range('E25').Value = '  +2.28%  '
$ws.Range('E26').Value = '  -0.01%  '
$c = $ws.Range('D27')
$c.Value = '''25.23'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +1.11%  '
$ws.Range('E28').Value = '  -2.87%  '
$ws.Range('E29').Value = '  +4.66%  '
$c = $ws.Range('D30')
$c.Value = '''33.14'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +5.12%  '
$ws.Range('B31').Value = 'Celestia'
$ws.Range('C31').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Range('D31')
$c.Value = '''19.52'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +11.09%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D32')
$c.Value = '''0.121'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +14.84%  '
$ws.Range('E33').Value = '  +2.35%  '
$ws.Range('E34').Value = '  +0.26%  '
$c = $ws.Range('D35')
$c.Value = '''0.0765'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +3.66%  '
$ws.Range('E36').Value = '  +2.24%  '
$c = $ws.Range('D37')
$c.Value = '''4.49'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +2.57%  '
$ws.Range('E38').Value = '  +3.43%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D39')
$c.Value = '''126.19'
$c.Style = "Normal"
$ws.Range('E39').Value = '  +9.60%  '
$ws.Range('B40').Value = 'WEMIXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D40')
$c.Value = '''2.30'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -0.53%  '
$ws.Range('E41').Value = '  +0.67%  '
$c = $ws.Range('D42')
$c.Value = '''21.94'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -1.98%  '
$ws.Range('E43').Value = '  +3.57%  '
$c = $ws.Range('D44')
$c.Value = '''1.946.37'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -0.36%  '
$c = $ws.Range('D45')
$c.Value = '''2.18'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +2.21%  '
$ws.Range('E46').Value = '  +6.91%  '
$c = $ws.Range('D47')
$c.Value = '''9.55'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +4.14%  '
$c = $ws.Range('D48')
$c.Value = '''1.68'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +9.88%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c = $ws.Range('D49')
$c.Value = '''2.673.61'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +2.62%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$c = $ws.Range('D50')
$c.Value = '''53.45'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +1.93%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Range('D51')
$c.Value = '''73.94'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +2.10%  '
